$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the two new rows of "persona" user-story content.
# Set A13 before A12 so the shared-string table gets the same
# insertion order as the target workbook (index 10 = "2 personas",
# index 11 = "invent users; ...").
$ws.Range("A13").Value = "2 personas"
$ws.Range("A12").Value = "invent users; give name,age, bg; backstory; technologically savvy or not; purpose of visit"

# Highlight the new cells with a solid yellow fill.
$ws.Range("A12:B12").Interior.Color = 65535
$ws.Range("A13").Interior.Color = 65535

# Update the active selection to match where the user ended up editing.
$ws.Range("B13").Select() | Out-Null
